$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.681.25"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +6.71%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.943.04"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +3.08%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.89"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.54%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.694"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.25%  "

$ws.Range("E7").Value = "  -0.12%  "

$ws.Range("E8").Value = "  +13.09%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.384"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +8.43%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "58.90"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +7.12%  "

$ws.Range("E11").Value = "  +4.16%  "

$ws.Range("E12").Value = "  +2.30%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.86"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +14.19%  "

$ws.Range("E14").Value = "  +9.37%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.225.14"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.99%  "

$ws.Range("E16").Value = "  +4.26%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.937.47"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.35%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.778.85"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +6.96%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "75.75"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +3.27%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0867"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +5.36%  "

$ws.Range("E21").Value = "  +7.36%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "254.68"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +4.13%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.24"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.67%  "

$ws.Range("E24").Value = "  +0.08%  "

$ws.Range("E25").Value = "  -6.12%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.83"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.55%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.13"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.85%  "

$ws.Range("E28").Value = "  +5.09%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.04"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +4.26%  "

$ws.Range("E30").Value = "  +1.78%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.63"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +8.25%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0620"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +4.37%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0924"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +28.73%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.37"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +3.74%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "19.69"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +43.08%  "

$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.89"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.96%  "

$ws.Range("B37").Value = "BinanceUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.06%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.902"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +6.08%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.45"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.48%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.01"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +4.08%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "105.76"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +7.98%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0228"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +2.94%  "

$ws.Range("E43").Value = "  +2.61%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.89"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +20.20%  "

$ws.Range("E45").Value = "  +4.37%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.359.83"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.74%  "

$ws.Range("E47").Value = "  +2.53%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0847"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +4.74%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.82"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +3.12%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.91"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +17.98%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.45"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.74%  "
